$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "State"
$ws.Cells.Item(1, 2).Value = "City"

# Data rows - State/City pairs (replaces prior Comunidad/Provincia Spanish data)
$data = @(
    @("California", "Los Angeles"),
    @("California", "San Diego"),
    @("California", "San Francisco"),
    @("California", "Sacramento"),
    @("California", "Death Valley"),
    @("Florida", "Jacksonville"),
    @("Florida", "Miami"),
    @("Florida", "Orlando"),
    @("Florida", "Tampa"),
    @("Illinois", "Cairo"),
    @("Illinois", "Chicago"),
    @("Illinois", "Rockford"),
    @("Illinois", "Springfield"),
    @("Nevada", "Las Vegas"),
    @("Nevada", "Reno"),
    @("Nevada", "Carson City"),
    @("Texas", "Dallas"),
    @("Texas", "Houston"),
    @("Texas", "El Paso"),
    @("Texas", "Brownsville")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Autofit column B to match the new, longer city names ("San Francisco")
$ws.Columns.Item(2).AutoFit() | Out-Null

# Move the active selection like the author's session ended
$ws.Range("F4").Select() | Out-Null
